$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: KARIM JALAL -> ANAS MASTI
$ws.Range("A2").Value = "ANAS MASTI"
$ws.Range("B2").Value = "BK747A53"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "845777567757575888678487"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "CIH"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "456/CASA"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2000

# Row 3: NABIL KAMAL -> JAMAL JAMAL
$ws.Range("A3").Value = "JAMAL JAMAL"
$ws.Range("B3").Value = "GT744635"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "474747446474747474747474"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "CIH"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "456/CASA"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1000

# Row 4: KHADIJA LALA -> blank (single space) row with new totals
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("I4").Value = 3000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3000

# Remove old row 5 (former totals row) entirely, shrinking the used range to A1:K4
$ws.Rows.Item(5).Delete()
